# Slide 3 ("cascade figures") - update german/bank/compas data labels
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$grp = $s.Shapes.Item(1)           # "Group 1" containing the cascade figure shapes

# TextBox 16: "3" -> "4"
$grp.GroupItems.Item(8).TextFrame.TextRange.Text = "4"

# TextBox 17: reposition/resize and "9" -> "12"
$tb17 = $grp.GroupItems.Item(9)
$tb17.Left = 433.5275
$tb17.Top = 219.316
$tb17.Width = 45.8809
$tb17.Height = 36.3516
$tb17.TextFrame.TextRange.Text = "12"

# TextBox 18: "6" -> "8"
$grp.GroupItems.Item(10).TextFrame.TextRange.Text = "8"

# TextBox 20: "18" -> "24"
$grp.GroupItems.Item(12).TextFrame.TextRange.Text = "24"
